$wb = $excel.ActiveWorkbook

# --- Update the "Date" metadata value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-02-20T19:32:51+00:00"

# --- Add a new concept row ("TRBA" / Transcriptome Bioinformatic Analysis) ---
# to the Concepts sheet, matching the formatting of the existing rows.
$concepts = $wb.Worksheets.Item("Concepts")

# Copy formatting (not values) from the last existing data row (row 5)
# down into the new row 6 so D6 and the rest of the row pick up the
# correct cell style.
$concepts.Range("A5:D5").Copy()
$concepts.Range("A6:D6").PasteSpecial(-4122)

# Copy the "Level" value (shared string "1") from A5 into A6 so it
# reuses the same shared-string entry instead of becoming a numeric literal.
$concepts.Range("A5").Copy()
$concepts.Range("A6").PasteSpecial(-4163)

# Fill in the new Code / Display values.
$concepts.Range("B6").Value = "TRBA"
$concepts.Range("C6").Value = "Transcriptome Bioinformatic Analysis"
